$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while forcing a genuine text cell
# (avoids Excel's automatic "False"/"True"/numeric-literal type coercion
# that would otherwise turn the text into a boolean/number cell).
function Set-TextValue($cellRef, $text) {
    $helper = $ws.Range("Z1000")
    $helper.Formula = '=T("' + $text + '")'
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helper.ClearContents()
}

# New variable row (row 6) matching the existing table's layout:
# Id | Name | Source | Type | MeasuringUnit | Comment | Alarm | AlarmGroup | AlarmText | AlarmLimitMin | AlarmLimitMax
$ws.Range("A6").Value = 5
Set-TextValue "B6" "TI_1"
Set-TextValue "D6" "REAL"
Set-TextValue "E6" "%"
Set-TextValue "F6" "Tank T1 - Level"
Set-TextValue "C6" "DB1.DBD6"
Set-TextValue "G6" "False"
Set-TextValue "H6" "None"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0

# Match the same row formatting (vertical-centered) used by the other data rows
$ws.Range("A6:K6").VerticalAlignment = -4108

# Update the active selection to match the saved workbook state
$ws.Range("H9").Select() | Out-Null
